$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.3
$ws.Range("H3").Value = 2.9
$ws.Range("I3").Value = 3.6
$ws.Range("K3").Value = 1.8
$ws.Range("L3").Value = 4.5
$ws.Range("Q3").Value = 3.1
$ws.Range("R3").Value = 1.36
$ws.Range("S3").Value = 5.4
$ws.Range("W3").Value = 1.73
$ws.Range("X3").Value = 2
$ws.Range("Y3").Value = 2.5
$ws.Range("Z3").Value = 1.5
$ws.Range("AA3").Value = 5
$ws.Range("AB3").Value = 9
$ws.Range("AG3").Value = 5
$ws.Range("AI3").Value = 23
$ws.Range("AJ3").Value = 101
$ws.Range("AN3").Value = 15

# Row 4
$ws.Range("G4").Value = 2.7
$ws.Range("I4").Value = 2.8
$ws.Range("J4").Value = 3.4
$ws.Range("L4").Value = 3.5
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.6
$ws.Range("S4").Value = 3.65
$ws.Range("T4").Value = 1.29
$ws.Range("U4").Value = 4.33
$ws.Range("V4").Value = 1.2
$ws.Range("Y4").Value = 1.95
$ws.Range("Z4").Value = 1.8
$ws.Range("AA4").Value = 7.5
$ws.Range("AE4").Value = 23
$ws.Range("AG4").Value = 7.5
$ws.Range("AH4").Value = 6
$ws.Range("AI4").Value = 15
$ws.Range("AJ4").Value = 51
$ws.Range("AK4").Value = 401
$ws.Range("AR4").Value = 1.78
$ws.Range("AS4").Value = 2.1

# Row 5
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 3
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 1.17
$ws.Range("N5").Value = 5
$ws.Range("AG5").Value = 5
$ws.Range("AH5").Value = 6
$ws.Range("AS5").Value = 1.63

# Row 6
$ws.Range("I6").Value = 4.2
$ws.Range("J6").Value = 2.88
$ws.Range("O6").Value = 1.62
$ws.Range("P6").Value = 2.2
$ws.Range("Q6").Value = 2.88
$ws.Range("R6").Value = 1.4
$ws.Range("S6").Value = 4.9
$ws.Range("T6").Value = 1.18
$ws.Range("AB6").Value = 8
$ws.Range("AI6").Value = 21

# Row 7
$ws.Range("G7").Value = 6.5
$ws.Range("I7").Value = 1.57
$ws.Range("Y7").Value = 2.2
$ws.Range("Z7").Value = 1.62
$ws.Range("AC7").Value = 21
$ws.Range("AJ7").Value = 81

# Row 9
$ws.Range("G9").Value = 1.4
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 9.5
$ws.Range("J9").Value = 2
$ws.Range("L9").Value = 9
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("AA9").Value = 4.75
$ws.Range("AD9").Value = 8.5
$ws.Range("AH9").Value = 8.5
$ws.Range("AJ9").Value = 126
$ws.Range("AL9").Value = 17
$ws.Range("AO9").Value = 126

# Row 10
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 2.45
$ws.Range("J10").Value = 4
$ws.Range("AO10").Value = 23
$ws.Range("AP10").Value = 23

# Row 11
$ws.Range("G11").Value = 3.05
$ws.Range("H11").Value = 2.35
$ws.Range("I11").Value = 3.05
$ws.Range("J11").Value = 3.95
$ws.Range("K11").Value = 1.65
$ws.Range("L11").Value = 3.95
$ws.Range("M11").Value = 1.21
$ws.Range("N11").Value = 3.95
$ws.Range("O11").Value = 1.8
$ws.Range("P11").Value = 1.91
$ws.Range("Q11").Value = 3.3
$ws.Range("R11").Value = 1.29
$ws.Range("U11").Value = 6.1
$ws.Range("W11").Value = 1.78
$ws.Range("X11").Value = 1.93
$ws.Range("Y11").Value = 2.35
$ws.Range("Z11").Value = 1.52
$ws.Range("AA11").Value = 5.8
$ws.Range("AB11").Value = 13.5
$ws.Range("AD11").Value = 45
$ws.Range("AF11").Value = 70
$ws.Range("AG11").Value = 3.95
$ws.Range("AH11").Value = 5.1
$ws.Range("AI11").Value = 20
$ws.Range("AM11").Value = 13.5
$ws.Range("AN11").Value = 12.5
$ws.Range("AO11").Value = 45
$ws.Range("AP11").Value = 40

# Row 13
$ws.Range("G13").Value = 2.55
$ws.Range("I13").Value = 3.2
$ws.Range("J13").Value = 3.15
$ws.Range("L13").Value = 3.85
$ws.Range("O13").Value = 1.53
$ws.Range("P13").Value = 2.35
$ws.Range("Q13").Value = 2.55
$ws.Range("R13").Value = 1.45
$ws.Range("U13").Value = 4.55
$ws.Range("X13").Value = 2.32
$ws.Range("Y13").Value = 2.02
$ws.Range("AA13").Value = 6.3
$ws.Range("AC13").Value = 9.75
$ws.Range("AD13").Value = 29
$ws.Range("AE13").Value = 25
$ws.Range("AF13").Value = 40
$ws.Range("AJ13").Value = 100
$ws.Range("AL13").Value = 7
$ws.Range("AM13").Value = 15
$ws.Range("AN13").Value = 11.75
$ws.Range("AO13").Value = 45
$ws.Range("AP13").Value = 37

# Row 14
$ws.Range("G14").Value = 1.4
$ws.Range("H14").Value = 4.05
$ws.Range("I14").Value = 8.25
$ws.Range("J14").Value = 1.88
$ws.Range("K14").Value = 2.22
$ws.Range("N14").Value = 6.8
$ws.Range("O14").Value = 1.35
$ws.Range("P14").Value = 2.95
$ws.Range("Q14").Value = 2.02
$ws.Range("R14").Value = 1.72
$ws.Range("U14").Value = 3.4
$ws.Range("V14").Value = 1.27
$ws.Range("W14").Value = 1.42
$ws.Range("X14").Value = 2.67
$ws.Range("Y14").Value = 2.3
$ws.Range("Z14").Value = 1.55
$ws.Range("AB14").Value = 5.6
$ws.Range("AD14").Value = 8.5
$ws.Range("AF14").Value = 37
$ws.Range("AG14").Value = 6.8
$ws.Range("AH14").Value = 8.25
$ws.Range("AL14").Value = 16
$ws.Range("AN14").Value = 27

# Row 15
$ws.Range("AH15").Value = 7
$ws.Range("AI15").Value = 17
$ws.Range("AL15").Value = 12

# Row 16
$ws.Range("G16").Value = 1.7
$ws.Range("I16").Value = 5.5
$ws.Range("J16").Value = 2.4
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 1.1
$ws.Range("N16").Value = 7
$ws.Range("W16").Value = 1.53
$ws.Range("X16").Value = 2.38
$ws.Range("Y16").Value = 2.25
$ws.Range("Z16").Value = 1.57
$ws.Range("AB16").Value = 7
$ws.Range("AD16").Value = 13
$ws.Range("AF16").Value = 41
$ws.Range("AG16").Value = 7
$ws.Range("AR16").Value = 1.8
$ws.Range("AS16").Value = 2.05

# Row 17
$ws.Range("M17").Value = 1.11
$ws.Range("N17").Value = 6.5
$ws.Range("O17").Value = 1.53
$ws.Range("P17").Value = 2.5
$ws.Range("Q17").Value = 2.75
$ws.Range("R17").Value = 1.44
$ws.Range("U17").Value = 5.5
$ws.Range("V17").Value = 1.14
$ws.Range("AR17").Value = 2.03
$ws.Range("AS17").Value = 1.83

# Row 18
$ws.Range("G18").Value = 1.6
$ws.Range("N18").Value = 8.5
$ws.Range("AD18").Value = 11
$ws.Range("AH18").Value = 7.5
$ws.Range("AI18").Value = 21
$ws.Range("AM18").Value = 29
$ws.Range("AN18").Value = 19
$ws.Range("AO18").Value = 67
$ws.Range("AP18").Value = 51

# Row 20
$ws.Range("L20").Value = 3.75
$ws.Range("Q20").Value = 2
$ws.Range("U20").Value = 3.75
$ws.Range("V20").Value = 1.29

# Row 21
$ws.Range("G21").Value = 1.95
$ws.Range("H21").Value = 3.3
$ws.Range("I21").Value = 4
$ws.Range("J21").Value = 2.63
$ws.Range("AC21").Value = 9

# Row 23
$ws.Range("G23").Value = 2
$ws.Range("I23").Value = 4.1
$ws.Range("J23").Value = 2.88
$ws.Range("Q23").Value = 2.5
$ws.Range("R23").Value = 1.5
$ws.Range("W23").Value = 1.62
$ws.Range("X23").Value = 2.2
$ws.Range("AC23").Value = 10
$ws.Range("AE23").Value = 21
$ws.Range("AG23").Value = 6
$ws.Range("AH23").Value = 6.5
$ws.Range("AL23").Value = 8.5
$ws.Range("AM23").Value = 19
$ws.Range("AN23").Value = 15
$ws.Range("AO23").Value = 41
$ws.Range("AR23").Value = 1.93
$ws.Range("AS23").Value = 1.88

# Row 24
$ws.Range("G24").Value = 1.7
$ws.Range("H24").Value = 3.7
$ws.Range("I24").Value = 4.75
$ws.Range("J24").Value = 2.3
$ws.Range("Q24").Value = 1.9
$ws.Range("R24").Value = 1.9
$ws.Range("AG24").Value = 11

# Row 25
$ws.Range("O25").Value = 1.36
$ws.Range("P25").Value = 3
$ws.Range("Q25").Value = 2.15
$ws.Range("R25").Value = 1.67
